# edit.ps1 - applies the "added new monitors and fixed filters" commit to monitors.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Data fix: H37 persistence value 1.5 -> 1
# ---------------------------------------------------------------------------
$ws.Cells.Item(37, 8).Value2 = 1

# ---------------------------------------------------------------------------
# 2) Add five new monitor rows (70-74).
#    The order in which brand-new text values are written matters because it
#    determines the order new entries are appended to the shared string
#    table, so we write the "new text" cells in the same sequence the
#    original author did, interleaved with the rest of each row.
# ---------------------------------------------------------------------------

# --- Row 70: LG 32GQ850 ---
$ws.Cells.Item(70, 1).Value2  = "LG 32GQ850"          # A70 (new string)
$ws.Cells.Item(70, 2).Value2  = 7.7
$ws.Cells.Item(70, 3).Value2  = 6
$ws.Cells.Item(70, 4).Value2  = 2
$ws.Cells.Item(70, 5).Value2  = 3
$ws.Cells.Item(70, 6).Value2  = 2
$ws.Cells.Item(70, 7).Value2  = 5
$ws.Cells.Item(70, 8).Value2  = 3
$ws.Cells.Item(70, 9).Value2  = "2560x1440"
$ws.Cells.Item(70, 10).Value2 = "260hz"               # J70 (new string)

# --- Row 71: Gigabyte FI32Q-X ---
$ws.Cells.Item(71, 10).Value2 = "270hz"               # J71 (new string, before A71)
$ws.Cells.Item(71, 1).Value2  = "Gigabyte FI32Q-X"    # A71 (new string)
$ws.Cells.Item(71, 2).Value2  = 7.8
$ws.Cells.Item(71, 3).Value2  = 6
$ws.Cells.Item(71, 4).Value2  = 2
$ws.Cells.Item(71, 5).Value2  = 3
$ws.Cells.Item(71, 6).Value2  = 2
$ws.Cells.Item(71, 7).Value2  = 5
$ws.Cells.Item(71, 8).Value2  = 3
$ws.Cells.Item(71, 9).Value2  = "2560x1440"

# --- finish off the review links for rows 70 & 71 ---
$ws.Cells.Item(70, 20).Value2 = "Bijan Jamshidi,https://www.youtube.com/watch?v=-uv7io23Dsg"                    # T70 (new string)
$ws.Cells.Item(71, 20).Value2 = "TFTCentral,https://tftcentral.co.uk/reviews/gigabyte-aorus-fi32q-x"            # T71 (new string)

# --- fill remaining shared (re-used) columns for rows 70 & 71 ---
foreach ($r in 70,71) {
    $ws.Cells.Item($r, 11).Value2 = "IPS"
    $ws.Cells.Item($r, 12).Value2 = "32"""
    $ws.Cells.Item($r, 13).Value2 = 700
    $ws.Cells.Item($r, 14).Value2 = 5700
    $ws.Cells.Item($r, 15).Value2 = "No"
    $ws.Cells.Item($r, 16).Value2 = "No"
    $ws.Cells.Item($r, 17).Value2 = "No"
    $ws.Cells.Item($r, 18).Value2 = "Wide"
    $ws.Cells.Item($r, 19).Value2 = "no"
}

# --- Row 72: Acer XB323U-GX ---
$ws.Cells.Item(72, 1).Value2  = "Acer XB323U-GX"      # A72 (new string)
$ws.Cells.Item(72, 2).Value2  = 7.8
$ws.Cells.Item(72, 3).Value2  = 6
$ws.Cells.Item(72, 4).Value2  = 2
$ws.Cells.Item(72, 5).Value2  = 3
$ws.Cells.Item(72, 6).Value2  = 2
$ws.Cells.Item(72, 7).Value2  = 5
$ws.Cells.Item(72, 8).Value2  = 3
$ws.Cells.Item(72, 9).Value2  = "2560x1440"
$ws.Cells.Item(72, 10).Value2 = "270hz"
$ws.Cells.Item(72, 11).Value2 = "IPS"
$ws.Cells.Item(72, 12).Value2 = "32"""
$ws.Cells.Item(72, 13).Value2 = 700
$ws.Cells.Item(72, 14).Value2 = 5700
$ws.Cells.Item(72, 15).Value2 = "No"
$ws.Cells.Item(72, 16).Value2 = "No"
$ws.Cells.Item(72, 17).Value2 = "No"
$ws.Cells.Item(72, 18).Value2 = "Wide"
$ws.Cells.Item(72, 19).Value2 = "Same tuning and panel as 32GQ850"   # S72 (new string)
$ws.Cells.Item(72, 20).Value2 = "no"

# --- Row 73: LG 42C2 ---
$ws.Cells.Item(73, 1).Value2  = "LG 42C2"             # A73 (new string)
$ws.Cells.Item(73, 2).Value2  = 5
$ws.Cells.Item(73, 3).Value2  = 9.5
$ws.Cells.Item(73, 4).Value2  = 10
$ws.Cells.Item(73, 5).Value2  = 7
$ws.Cells.Item(73, 6).Value2  = 3.5
$ws.Cells.Item(73, 7).Value2  = 8
$ws.Cells.Item(73, 8).Value2  = 1.5
$ws.Cells.Item(73, 9).Value2  = "3840x2160"
$ws.Cells.Item(73, 10).Value2 = "120hz"
$ws.Cells.Item(73, 11).Value2 = "W-OLED"
$ws.Cells.Item(73, 12).Value2 = "42"""                # L73 (new string)
$ws.Cells.Item(73, 13).Value2 = 1000
$ws.Cells.Item(73, 14).Value2 = "ps5"
$ws.Cells.Item(73, 15).Value2 = "No"
$ws.Cells.Item(73, 16).Value2 = "No"
$ws.Cells.Item(73, 17).Value2 = "Yes"
$ws.Cells.Item(73, 18).Value2 = "Wide"
$ws.Cells.Item(73, 20).Value2 = "RTINGS,https://www.rtings.com/monitor/reviews/lg/42-c2-oled;Hardware Unboxed,https://www.youtube.com/watch?v=jRzGvkqSNaI"  # T73 (new string)

# --- Row 74: Sony Inzone M9 ---
$ws.Cells.Item(74, 1).Value2  = "Sony Inzone M9"      # A74 (new string)
$ws.Cells.Item(74, 2).Value2  = 5.8
$ws.Cells.Item(74, 3).Value2  = 6.2
$ws.Cells.Item(74, 4).Value2  = 4
$ws.Cells.Item(74, 5).Value2  = 6.1
$ws.Cells.Item(74, 6).Value2  = 4
$ws.Cells.Item(74, 7).Value2  = 8
$ws.Cells.Item(74, 8).Value2  = 8
$ws.Cells.Item(74, 9).Value2  = "3840x2160"
$ws.Cells.Item(74, 10).Value2 = "144hz"
$ws.Cells.Item(74, 11).Value2 = "IPS FALD 96 zones"   # K74 (new string)
$ws.Cells.Item(74, 12).Value2 = "27"""
$ws.Cells.Item(74, 13).Value2 = 900
$ws.Cells.Item(74, 14).Value2 = "ps5"
$ws.Cells.Item(74, 15).Value2 = "No"
$ws.Cells.Item(74, 16).Value2 = "No"
$ws.Cells.Item(74, 17).Value2 = "No"
$ws.Cells.Item(74, 18).Value2 = "Wide"
$ws.Cells.Item(74, 19).Value2 = "Very entry-level HDR and quite overpriced"   # S74 (new string)
$ws.Cells.Item(74, 20).Value2 = "RTINGS,https://www.rtings.com/monitor/reviews/sony/inzone-m9;Hardware Unboxed,https://www.youtube.com/watch?v=GNF2YMuITr0"  # T74 (new string)

# --- last new string: the "special" note that was added to row 73 afterwards ---
$ws.Cells.Item(73, 19).Value2 = "Burn-in risk + glossy"   # S73 (new string, appended last)

# ---------------------------------------------------------------------------
# 3) Column K (11) gets a custom width
# ---------------------------------------------------------------------------
$ws.Columns.Item(11).ColumnWidth = 14.33

# ---------------------------------------------------------------------------
# 4) Stray formatted (but empty) cell E78 with a top border, left behind by
#    the author below the data block.
# ---------------------------------------------------------------------------
$e78 = $ws.Range("E78")
$e78.Borders.Item(8).LineStyle = 1
$e78.Borders.Item(8).Weight = 2

# ---------------------------------------------------------------------------
# 5) Update the view: scroll position and active selection
# ---------------------------------------------------------------------------
$win = $wb.Windows.Item(1)
$win.ScrollRow = 23
$win.ScrollColumn = 1
$ws.Range("S73").Select() | Out-Null
